$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the value to be stored as text (matches the inlineStr/shared-string
    # cell type used throughout this sheet), even when the value looks like a
    # number (e.g. "1.001", "0.6309"). A leading apostrophe marks it as text
    # input (quotePrefix), then we reset the cell style back to Normal so we
    # don't leave a stray formatting difference behind.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '29.449.18'
Set-TextValue $ws.Range('E2') '  +0.21%  '
Set-TextValue $ws.Range('D3') '1.852.79'
Set-TextValue $ws.Range('E3') '  +0.31%  '
Set-TextValue $ws.Range('D5') '240.83'
Set-TextValue $ws.Range('E5') '  +0.86%  '
Set-TextValue $ws.Range('D6') '0.6309'
Set-TextValue $ws.Range('E6') '  +0.13%  '
Set-TextValue $ws.Range('E7') '  -0.01%  '
Set-TextValue $ws.Range('D8') '0.07676'
Set-TextValue $ws.Range('E8') '  +1.74%  '
Set-TextValue $ws.Range('D9') '0.2944'
Set-TextValue $ws.Range('E9') '  -0.07%  '
Set-TextValue $ws.Range('D10') '24.69'
Set-TextValue $ws.Range('E10') '  +0.16%  '
Set-TextValue $ws.Range('D11') '0.07756'
Set-TextValue $ws.Range('E11') '  +0.73%  '
Set-TextValue $ws.Range('D12') '1.857.57'
Set-TextValue $ws.Range('E12') '  +0.33%  '
Set-TextValue $ws.Range('D13') '5.032'
Set-TextValue $ws.Range('E13') '  +0.89%  '
Set-TextValue $ws.Range('D14') '0.00001091'
Set-TextValue $ws.Range('E14') '  +7.27%  '
Set-TextValue $ws.Range('D15') '0.6818'
Set-TextValue $ws.Range('E15') '  +0.27%  '
Set-TextValue $ws.Range('D16') '83.75'
Set-TextValue $ws.Range('E16') '  +0.69%  '
Set-TextValue $ws.Range('D17') '2.099.89'
Set-TextValue $ws.Range('E17') '  -0.52%  '
Set-TextValue $ws.Range('D18') '6.170'
Set-TextValue $ws.Range('E18') '  +0.60%  '
Set-TextValue $ws.Range('D19') '29.478.83'
Set-TextValue $ws.Range('E19') '  +0.17%  '
Set-TextValue $ws.Range('D20') '230.00'
Set-TextValue $ws.Range('E20') '  +0.91%  '
Set-TextValue $ws.Range('D21') '12.48'
Set-TextValue $ws.Range('E21') '  +0.03%  '
Set-TextValue $ws.Range('D22') '1.001'
Set-TextValue $ws.Range('E22') '  -0.02%  '
Set-TextValue $ws.Range('D23') '7.461'
Set-TextValue $ws.Range('E23') '  -0.51%  '
Set-TextValue $ws.Range('E24') '  -0.09%  '
Set-TextValue $ws.Range('D25') '156.97'
Set-TextValue $ws.Range('E25') '  +0.04%  '
Set-TextValue $ws.Range('E26') '  -0.57%  '
Set-TextValue $ws.Range('D27') '8.418'
Set-TextValue $ws.Range('E27') '  +0.58%  '
Set-TextValue $ws.Range('D28') '17.75'
Set-TextValue $ws.Range('E28') '  +0.71%  '
Set-TextValue $ws.Range('D29') '1.322'
Set-TextValue $ws.Range('E29') '  +3.85%  '
Set-TextValue $ws.Range('D30') '1.469'
Set-TextValue $ws.Range('E30') '  +0.46%  '
Set-TextValue $ws.Range('D31') '0.05689'
Set-TextValue $ws.Range('E31') '  +0.38%  '
Set-TextValue $ws.Range('D32') '4.132'
Set-TextValue $ws.Range('E32') '  +0.19%  '
Set-TextValue $ws.Range('E33') '  +0.55%  '
Set-TextValue $ws.Range('D34') '1.853'
Set-TextValue $ws.Range('E34') '  +0.80%  '
Set-TextValue $ws.Range('D35') '1.164'
Set-TextValue $ws.Range('E35') '  +0.67%  '
Set-TextValue $ws.Range('D36') '0.7054'
Set-TextValue $ws.Range('E36') '  -1.38%  '
Set-TextValue $ws.Range('D38') '2.782'
Set-TextValue $ws.Range('E38') '  +0.23%  '
Set-TextValue $ws.Range('D39') '0.01794'
Set-TextValue $ws.Range('E39') '  -0.59%  '
Set-TextValue $ws.Range('D40') '1.219.76'
Set-TextValue $ws.Range('E40') '  -2.26%  '
Set-TextValue $ws.Range('D41') '6.541'
Set-TextValue $ws.Range('E41') '  +5.65%  '
Set-TextValue $ws.Range('D42') '0.9076'
Set-TextValue $ws.Range('E42') '  +0.09%  '
Set-TextValue $ws.Range('E43') '  +0.02%  '
Set-TextValue $ws.Range('B44') 'Quant'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D44') '101.76'
Set-TextValue $ws.Range('E44') '  +0.16%  '
Set-TextValue $ws.Range('B45') 'Aave'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D45') '66.49'
Set-TextValue $ws.Range('E45') '  +0.35%  '
Set-TextValue $ws.Range('B46') 'BabyDogeCoin'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D46') '0.00000000120'
Set-TextValue $ws.Range('E46') '  +1.49%  '
Set-TextValue $ws.Range('B47') 'Aptos'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D47') '7.126'
Set-TextValue $ws.Range('E47') '  +0.78%  '
Set-TextValue $ws.Range('B48') 'TheSandbox'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D48') '0.4020'
Set-TextValue $ws.Range('E48') '  +0.16%  '
Set-TextValue $ws.Range('B49') 'EnergySwap'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D49') '9.031'
Set-TextValue $ws.Range('E49') '  -0.08%  '
Set-TextValue $ws.Range('B50') 'RenderToken'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D50') '1.686'
Set-TextValue $ws.Range('E50') '  -0.32%  '
Set-TextValue $ws.Range('B51') 'Algorand'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D51') '0.1135'
Set-TextValue $ws.Range('E51') '  +1.61%  '
